# Video merger spreadsheet: the "audio track" column (D) was pointed at
# videos/audio1.mp3 for every data row; update it to videos/audio0.mp3
# (output columns A-C, E are unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D1").Value = "videos/audio0.mp3"
$ws.Range("D2").Value = "videos/audio0.mp3"
$ws.Range("D3").Value = "videos/audio0.mp3"

# Leave the cursor where the author last left it when saving.
$ws.Range("C3").Select()
